# Applies updated test-run data (durations, timestamps) to existing rows 2-11,
# and appends new rows 12-14 for additional test results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-11: duration (E) and timestamp (F) ---

$ws.Range("E2").Value = 43.60830139997415
$ws.Range("F2").Value = "2024-04-14T00:12:38"

$ws.Range("E3").Value = 0.4832969999988563
$ws.Range("F3").Value = "2024-04-14T00:12:38"

$ws.Range("E4").Value = 43.53136759999325
$ws.Range("F4").Value = "2024-04-14T00:13:22"

$ws.Range("E5").Value = 1.498724299977766
$ws.Range("F5").Value = "2024-04-14T00:13:23"

$ws.Range("E6").Value = 1.463173900003312
$ws.Range("F6").Value = "2024-04-14T00:13:25"

$ws.Range("E7").Value = 42.97459470000467
$ws.Range("F7").Value = "2024-04-14T00:14:52"

$ws.Range("E8").Value = 0.4309544999850914
$ws.Range("F8").Value = "2024-04-14T00:14:52"

$ws.Range("E9").Value = 0.4897482999949716
$ws.Range("F9").Value = "2024-04-14T00:14:53"

$ws.Range("E10").Value = 0.9210826000198722
$ws.Range("F10").Value = "2024-04-14T00:14:54"

$ws.Range("E11").Value = 1.028996900015045
$ws.Range("F11").Value = "2024-04-14T00:14:55"

# --- Append new rows 12-14 ---

# Row 12
$ws.Range("A12").Value = "TestUsers"
$ws.Range("B12").Value = "test_required_field_name"
$ws.Range("C12").Value = "Test required field is not sent in request body (name)"
$ws.Range("D12").Value = "PASSED"
$ws.Range("E12").Value = 0.4375171000137925
$ws.Range("F12").Value = "2024-04-14T00:14:55"
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = "gorest_api\users\test_users.py"
$ws.Range("I12").Value = "negative"

# Row 13
$ws.Range("A13").Value = "TestUsers"
$ws.Range("B13").Value = "test_email_already_taken"
$ws.Range("C13").Value = "Test email address is already taken"
$ws.Range("D13").Value = "PASSED"
$ws.Range("E13").Value = 0.9908925000054296
$ws.Range("F13").Value = "2024-04-14T00:14:56"
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = "gorest_api\users\test_users.py"
$ws.Range("I13").Value = "negative"

# Row 14
$ws.Range("A14").Value = "TestUsers"
$ws.Range("B14").Value = "test_nonexistent_user"
$ws.Range("C14").Value = "Test trying to retrieve a user that does not exist"
$ws.Range("D14").Value = "PASSED"
$ws.Range("E14").Value = 0.435685900010867
$ws.Range("F14").Value = "2024-04-14T00:14:56"
$ws.Range("G14").Value = ""
$ws.Range("H14").Value = "gorest_api\users\test_users.py"
$ws.Range("I14").Value = "negative"
